$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Update "VALOR MORA" total (E11): 136000 -> 204000
# ------------------------------------------------------------------
$ws.Range("E11").Value = 204000

# ------------------------------------------------------------------
# 2) Update "Cant. Periodos" (F13): 2 -> 3
# ------------------------------------------------------------------
$ws.Range("F13").Value = 3

# ------------------------------------------------------------------
# 3) Insert a new data row (row 18) for period 2508, re-using the
#    formatting that the old closing row (17) had, and move the
#    "closing" border look down to the new last row.
# ------------------------------------------------------------------
$ws.Rows("18").Insert()

# Copy the old row 17 (still intact just below the inserted blank row)
# into new row 18 - formats first, then values - so row 18 gets the
# "bottom of table" border style that row 17 used to have.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# Row 17 is no longer the last row of the table, so it now takes on
# the "inner" border look that row 16 uses.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Fix up the "Periodo Mora" values so the three rows read, in
#    order, 2506 / 2507 / 2508.
# ------------------------------------------------------------------
$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2507"
$ws.Range("E18").Value = "2508"
